$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ===================================================================
# Sheet "results" (sheet1): insert two new columns -
#   "S*-unmerged"   (a timing column, right after "S*-MM0", like B..F)
#   "S*-unmergedND" (a boolean column, right after "S*-MM0ND", like
#                    the other ND columns)
# ===================================================================

# --- Header row 1: shift G1..K1 one column right into H1..L1, then
#     put the two brand-new headers into G1 and M1. Style (bold /
#     border / centered) is copied from the existing K1 header cell
#     onto the two brand-new cells (L1, M1).

$ws1.Range("K1").Copy()
$ws1.Range("L1:M1").PasteSpecial(-4122)

$ws1.Range("L1").Value = "S*-MM0ND"
$ws1.Range("K1").Value = "S*-MMND"
$ws1.Range("J1").Value = "S*-HSND"
$ws1.Range("I1").Value = "S*-BSND"
$ws1.Range("H1").Value = "var"
$ws1.Range("G1").Value = "S*-unmerged"
$ws1.Range("M1").Value = "S*-unmergedND"

# --- Row 2 (data): shift H2..K2 one column right into I2..L2, then
#     set the new G2 timing value and the new trailing M2 boolean.
#     (Every one of the shifted/new ND columns ends up $false here.)

$ws1.Range("G2").Value = 273
$ws1.Range("H2").Value = 0
$ws1.Range("I2").Value = $false
$ws1.Range("J2").Value = $false
$ws1.Range("K2").Value = $false
$ws1.Range("L2").Value = $false
$ws1.Range("M2").Value = $false

# ===================================================================
# Sheet "stats" (sheet2): a new "S*-unmerged" row is inserted right
# before each "Kruskal" summary row (one in the "run 0" block, one in
# the "Average" block) and every numeric value on the sheet is
# refreshed for the new run.
# ===================================================================

# Break the two existing vertical merges before moving data under them.
$ws2.Range("A2:A6").UnMerge()
$ws2.Range("A7:A11").UnMerge()

# ---- "run 0" block: rows 2-7 ----
$ws2.Cells.Item(2, 2).Value = "S*-BS"
$ws2.Cells.Item(2, 3).Value = 84
$ws2.Cells.Item(2, 4).Value = 0.0001854621805250645
$ws2.Cells.Item(2, 5).Value = 0.0851059821434319
$ws2.Cells.Item(2, 6).Value = 84
$ws2.Cells.Item(2, 7).Value = 0.005574297625571489
$ws2.Cells.Item(2, 8).Value = 0.01304262317717075
$ws2.Cells.Item(2, 9).Value = 0.01777624059468508
$ws2.Cells.Item(2, 10).Value = 0.03801995515823364
$ws2.Cells.Item(2, 11).Value = 0.001721014734357595

$ws2.Cells.Item(3, 2).Value = "S*-HS"
$ws2.Cells.Item(3, 3).Value = 84
$ws2.Cells.Item(3, 4).Value = 0.004277261439710855
$ws2.Cells.Item(3, 5).Value = 0.06250949809327722
$ws2.Cells.Item(3, 6).Value = 84
$ws2.Cells.Item(3, 7).Value = 0.005107311066240072
$ws2.Cells.Item(3, 8).Value = 0.01221129298210144
$ws2.Cells.Item(3, 9).Value = 0.01376722613349557
$ws2.Cells.Item(3, 10).Value = 0.02505912631750107
$ws2.Cells.Item(3, 11).Value = 0.001718407962471247

$ws2.Cells.Item(4, 2).Value = "S*-MM"
$ws2.Cells.Item(4, 3).Value = 84
$ws2.Cells.Item(4, 4).Value = 0.003701074980199337
$ws2.Cells.Item(4, 5).Value = 0.05916893575340509
$ws2.Cells.Item(4, 6).Value = 84
$ws2.Cells.Item(4, 7).Value = 0.004579616244882345
$ws2.Cells.Item(4, 8).Value = 0.01076601119711995
$ws2.Cells.Item(4, 9).Value = 0.01292823813855648
$ws2.Cells.Item(4, 10).Value = 0.02521110186353326
$ws2.Cells.Item(4, 11).Value = 0.001566965598613024

$ws2.Cells.Item(5, 2).Value = "S*-MM0"
$ws2.Cells.Item(5, 3).Value = 84
$ws2.Cells.Item(5, 4).Value = 0.0002787499688565731
$ws2.Cells.Item(5, 5).Value = 0.1410107812844217
$ws2.Cells.Item(5, 6).Value = 84
$ws2.Cells.Item(5, 7).Value = 0.005227842833846807
$ws2.Cells.Item(5, 8).Value = 0.00830181036144495
$ws2.Cells.Item(5, 9).Value = 0.02027443004772067
$ws2.Cells.Item(5, 10).Value = 0.04659218247979879
$ws2.Cells.Item(5, 11).Value = 0.001684343907982111

# Give the brand-new row 6 the same plain bordered style as row 5
# (columns A and B) before filling it in.
$ws2.Range("A5:B5").Copy()
$ws2.Range("A6:B6").PasteSpecial(-4122)

$ws2.Cells.Item(6, 2).Value = "S*-unmerged"
$ws2.Cells.Item(6, 3).Value = 238
$ws2.Cells.Item(6, 4).Value = 0.006793376989662647
$ws2.Cells.Item(6, 5).Value = 0.2588178510777652
$ws2.Cells.Item(6, 6).Value = 238
$ws2.Cells.Item(6, 7).Value = 0.0124337007291615
$ws2.Cells.Item(6, 8).Value = 0.0299899373203516
$ws2.Cells.Item(6, 9).Value = 0.1408478491939604
$ws2.Cells.Item(6, 10).Value = 0.05743196373805404
$ws2.Cells.Item(6, 11).Value = 0.005592283327132463

# Brand-new row 7 ("Kruskal" summary row, shifted down from 6).
$ws2.Range("A5:B5").Copy()
$ws2.Range("A7:B7").PasteSpecial(-4122)

$ws2.Cells.Item(7, 2).Value = "Kruskal"
$ws2.Cells.Item(7, 3).Value = 2500
$ws2.Cells.Item(7, 4).ClearContents()
$ws2.Cells.Item(7, 5).Value = 0.05454195069614798
$ws2.Cells.Item(7, 6).ClearContents()
$ws2.Cells.Item(7, 7).ClearContents()
$ws2.Cells.Item(7, 8).ClearContents()
$ws2.Cells.Item(7, 9).ClearContents()
$ws2.Cells.Item(7, 10).ClearContents()
$ws2.Cells.Item(7, 11).ClearContents()

# ---- "Average" block: rows 8-13 (same values as the "run 0" block) ----
$ws2.Cells.Item(8, 2).Value = "S*-BS"
$ws2.Cells.Item(8, 3).Value = 84
$ws2.Cells.Item(8, 4).Value = 0.0001854621805250645
$ws2.Cells.Item(8, 5).Value = 0.0851059821434319
$ws2.Cells.Item(8, 6).Value = 84
$ws2.Cells.Item(8, 7).Value = 0.005574297625571489
$ws2.Cells.Item(8, 8).Value = 0.01304262317717075
$ws2.Cells.Item(8, 9).Value = 0.01777624059468508
$ws2.Cells.Item(8, 10).Value = 0.03801995515823364
$ws2.Cells.Item(8, 11).Value = 0.001721014734357595

$ws2.Cells.Item(9, 2).Value = "S*-HS"
$ws2.Cells.Item(9, 3).Value = 84
$ws2.Cells.Item(9, 4).Value = 0.004277261439710855
$ws2.Cells.Item(9, 5).Value = 0.06250949809327722
$ws2.Cells.Item(9, 6).Value = 84
$ws2.Cells.Item(9, 7).Value = 0.005107311066240072
$ws2.Cells.Item(9, 8).Value = 0.01221129298210144
$ws2.Cells.Item(9, 9).Value = 0.01376722613349557
$ws2.Cells.Item(9, 10).Value = 0.02505912631750107
$ws2.Cells.Item(9, 11).Value = 0.001718407962471247

$ws2.Cells.Item(10, 2).Value = "S*-MM"
$ws2.Cells.Item(10, 3).Value = 84
$ws2.Cells.Item(10, 4).Value = 0.003701074980199337
$ws2.Cells.Item(10, 5).Value = 0.05916893575340509
$ws2.Cells.Item(10, 6).Value = 84
$ws2.Cells.Item(10, 7).Value = 0.004579616244882345
$ws2.Cells.Item(10, 8).Value = 0.01076601119711995
$ws2.Cells.Item(10, 9).Value = 0.01292823813855648
$ws2.Cells.Item(10, 10).Value = 0.02521110186353326
$ws2.Cells.Item(10, 11).Value = 0.001566965598613024

$ws2.Cells.Item(11, 2).Value = "S*-MM0"
$ws2.Cells.Item(11, 3).Value = 84
$ws2.Cells.Item(11, 4).Value = 0.0002787499688565731
$ws2.Cells.Item(11, 5).Value = 0.1410107812844217
$ws2.Cells.Item(11, 6).Value = 84
$ws2.Cells.Item(11, 7).Value = 0.005227842833846807
$ws2.Cells.Item(11, 8).Value = 0.00830181036144495
$ws2.Cells.Item(11, 9).Value = 0.02027443004772067
$ws2.Cells.Item(11, 10).Value = 0.04659218247979879
$ws2.Cells.Item(11, 11).Value = 0.001684343907982111

# Brand-new row 12.
$ws2.Range("A11:B11").Copy()
$ws2.Range("A12:B12").PasteSpecial(-4122)

$ws2.Cells.Item(12, 2).Value = "S*-unmerged"
$ws2.Cells.Item(12, 3).Value = 238
$ws2.Cells.Item(12, 4).Value = 0.006793376989662647
$ws2.Cells.Item(12, 5).Value = 0.2588178510777652
$ws2.Cells.Item(12, 6).Value = 238
$ws2.Cells.Item(12, 7).Value = 0.0124337007291615
$ws2.Cells.Item(12, 8).Value = 0.0299899373203516
$ws2.Cells.Item(12, 9).Value = 0.1408478491939604
$ws2.Cells.Item(12, 10).Value = 0.05743196373805404
$ws2.Cells.Item(12, 11).Value = 0.005592283327132463

# Brand-new row 13 ("Kruskal" summary row, shifted down from 11).
$ws2.Range("A11:B11").Copy()
$ws2.Range("A13:B13").PasteSpecial(-4122)

$ws2.Cells.Item(13, 2).Value = "Kruskal"
$ws2.Cells.Item(13, 3).Value = 2500
$ws2.Cells.Item(13, 4).ClearContents()
$ws2.Cells.Item(13, 5).Value = 0.05454195069614798
$ws2.Cells.Item(13, 6).ClearContents()
$ws2.Cells.Item(13, 7).ClearContents()
$ws2.Cells.Item(13, 8).ClearContents()
$ws2.Cells.Item(13, 9).ClearContents()
$ws2.Cells.Item(13, 10).ClearContents()
$ws2.Cells.Item(13, 11).ClearContents()

# A2 / A8 carry the merged group labels ("run 0" / "Average").
$ws2.Cells.Item(2, 1).Value = "run 0"
$ws2.Cells.Item(8, 1).Value = "Average"

# Re-merge the (now taller) groups.
$ws2.Range("A2:A7").Merge()
$ws2.Range("A8:A13").Merge()
